$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.911.73'
$ws.Range("E2").Value = '  +2.23%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.109.08'
$ws.Range("E3").Value = '  +0.89%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.31%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.37'
$ws.Range("E5").Value = '  +0.97%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '616.44'
$ws.Range("E6").Value = '  -0.18%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.09'
$ws.Range("E7").Value = '  -3.96%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.395'
$ws.Range("E8").Value = '  +8.91%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.108.85'
$ws.Range("E10").Value = '  +1.00%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.734'
$ws.Range("E11").Value = '  +0.17%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.202'
$ws.Range("E12").Value = '  +0.11%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000253'
$ws.Range("E13").Value = '  +3.50%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.55'
$ws.Range("E14").Value = '  +0.06%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.182.66'
$ws.Range("E15").Value = '  +2.65%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.51'
$ws.Range("E16").Value = '  +1.02%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.691.45'
$ws.Range("E17").Value = '  +1.24%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.147.57'
$ws.Range("E18").Value = '  +2.31%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.62'
$ws.Range("E19").Value = '  -0.80%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.78'
$ws.Range("E20").Value = '  +2.43%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.80'
$ws.Range("E21").Value = '  +0.97%  '

# Row 22
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '447.40'
$ws.Range("E22").Value = '  +2.52%  '

# Row 23
$ws.Range("B23").Value = 'PEPE'
$ws.Range("C23").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000203'
$ws.Range("E23").Value = '  -2.21%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.29'
$ws.Range("E24").Value = '  +3.80%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.64'
$ws.Range("E25").Value = '  +0.93%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.64'
$ws.Range("E26").Value = '  -0.94%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '80.02'
$ws.Range("E27").Value = '  -11.57%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.278.78'

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.140'
$ws.Range("E30").Value = '  +17.88%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.229'
$ws.Range("E31").Value = '  -6.39%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.168'
$ws.Range("E32").Value = '  -3.67%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.32'
$ws.Range("E33").Value = '  +2.45%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.171'
$ws.Range("E34").Value = '  +3.16%  '

# Row 35
$ws.Range("E35").Value = '  -0.14%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.00'
$ws.Range("E36").Value = '  +5.09%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.33'
$ws.Range("E37").Value = '  +0.84%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.08'
$ws.Range("E38").Value = '  -5.17%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.93'
$ws.Range("E39").Value = '  +1.99%  '

# Row 40
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.31'
$ws.Range("E40").Value = '  +2.71%  '

# Row 41
$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '481.54'
$ws.Range("E41").Value = '  -0.54%  '

# Row 42
$ws.Range("B42").Value = 'PolygonEcosystemToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.436'
$ws.Range("E42").Value = '  +4.78%  '

# Row 43
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.47'
$ws.Range("E43").Value = '  -0.93%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.21'
$ws.Range("E44").Value = '  +0.29%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '158.85'
$ws.Range("E46").Value = '  +2.81%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.92'
$ws.Range("E47").Value = '  +2.09%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.698'
$ws.Range("E48").Value = '  +2.43%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.36'
$ws.Range("E49").Value = '  +2.26%  '

# Row 50
$ws.Range("E50").Value = '  +7.12%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.05'
$ws.Range("E51").Value = '  -0.15%  '
